$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated price / volume values scraped for this run.
# Each value is written with a leading apostrophe so Excel stores it as
# literal text (matching the original inlineStr cells) instead of parsing
# strings such as "26.973.73" or "0.4680" as numbers/dates, and the style
# is reset to "Normal" afterwards so no stray quote-prefix formatting is
# introduced.
$updates = @{
    "D2" = "26.973.73"
    "E2" = "  +0.29%  "
    "D3" = "1.820.20"
    "E3" = "  +0.22%  "
    "E4" = "  +0.21%  "
    "D5" = "310.96"
    "E5" = "  +0.23%  "
    "E6" = "  +0.19%  "
    "D7" = "0.4680"
    "E7" = "  +0.81%  "
    "D8" = "0.3667"
    "E8" = "  -0.93%  "
    "D9" = "0.07355"
    "E9" = "  +0.04%  "
    "E10" = "  +0.40%  "
    "D11" = "20.28"
    "E11" = "  -0.75%  "
    "D12" = "1.836.00"
    "E12" = "  +0.32%  "
    "D13" = "5.420"
    "E13" = "  +1.29%  "
    "D14" = "0.07170"
    "E14" = "  +0.96%  "
    "D15" = "6.513"
    "E15" = "  +0.11%  "
    "D16" = "91.62"
    "E16" = "  +0.15%  "
    "E17" = "  +0.38%  "
    "D18" = "0.000008749"
    "E18" = "  +0.43%  "
    "D20" = "14.67"
    "E20" = "  -0.25%  "
    "D21" = "26.998.51"
    "E21" = "  +0.23%  "
    "D22" = "5.290"
    "E22" = "  -0.71%  "
    "E23" = "  +0.49%  "
    "D24" = "2.045.56"
    "E24" = "  -0.55%  "
    "D25" = "1.891"
    "E25" = "  -0.51%  "
    "D26" = "150.82"
    "E26" = "  -0.71%  "
    "E27" = "  -0.01%  "
    "D28" = "2.144"
    "E28" = "  +0.21%  "
    "D29" = "5.233"
    "E29" = "  -1.29%  "
    "D30" = "116.63"
    "E30" = "  +1.09%  "
    "D31" = "0.08884"
    "E31" = "  -0.02%  "
    "E32" = "  -0.60%  "
    "D33" = "1.161"
    "E33" = "  +0.75%  "
    "E34" = "  +0.99%  "
    "D35" = "2.942"
    "E36" = "  +0.24%  "
    "D37" = "1.095"
    "E37" = "  +0.01%  "
    "D38" = "0.05312"
    "E38" = "  +0.92%  "
    "E39" = "  -0.48%  "
    "D40" = "2.980"
    "E40" = "  +1.48%  "
    "D41" = "2.372"
    "E41" = "  -0.44%  "
    "D42" = "7.186"
    "E42" = "  -0.64%  "
    "D43" = "0.5296"
    "E43" = "  -0.67%  "
    "D44" = "0.1652"
    "E44" = "  -0.61%  "
    "D45" = "8.461"
    "E45" = "  +0.15%  "
    "D46" = "0.4896"
    "E46" = "  -0.97%  "
    "D47" = "10.51"
    "E47" = "  +1.69%  "
    "E48" = "  +0.21%  "
    "D49" = "1.662"
    "E49" = "  -0.87%  "
    "D50" = "103.04"
    "E50" = "  -0.17%  "
    "D51" = "0.06296"
    "E51" = "  +0.19%  "
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $updates[$cellRef]
    $cell.Style = "Normal"
}

Write-Host "Updated $($updates.Count) cells"
